# Insert a new data row at row 209 (pushing existing rows 209-278 down to
# 210-279) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 209..278 down by one row (like right-clicking row 209's header
# and choosing "Insert").
$ws.Range("209:209").Insert()

# Fill in the newly inserted row 209 with the new record's data.
$ws.Range("A209").Value = 3
$ws.Range("B209").Value = "Femacal de La Calera"
$ws.Range("C209").Value = "Coquimbo"
$ws.Range("D209").Value = 44559
$ws.Range("E209").Value = 5
$ws.Range("F209").Value = 100112031
$ws.Range("G209").Value = "Poroto verde"
$ws.Range("H209").Value = "Magnum"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 87
$ws.Range("K209").Value = 23000
$ws.Range("L209").Value = 24000
$ws.Range("M209").Value = 23460
$ws.Range("N209").Value = '$/malla 25 kilos'
$ws.Range("O209").Value = "Provincia de Quillota"
$ws.Range("P209").Value = 938
$ws.Range("Q209").Value = 25
$ws.Range("R209").Value = "Hortaliza"
